$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.447252
$ws.Range("H2").Value = 61.341756
$ws.Range("I2").Value = 0.8699145605694745
$ws.Range("J2").Value = 0.8770588936480435
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.030493
$ws.Range("N2").Value = 0.091479
$ws.Range("O2").Value = 0.001018637778668347
$ws.Range("P2").Value = 0.001021037349570579
$ws.Range("Q2").Value = 0.623498055236
$ws.Range("R2").Value = 5.611482497124
$ws.Range("S2").Value = 0.0008861278356097409
$ws.Range("T2").Value = 0.0008955098881877029

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.447252
$ws.Range("H3").Value = 61.341756
$ws.Range("I3").Value = 0.8699145605694745
$ws.Range("J3").Value = 0.8770588936480435
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 29.65321533333333
$ws.Range("N3").Value = 88.95964599999999
$ws.Range("O3").Value = 0.990584245483253
$ws.Range("P3").Value = 0.9929177316168408
$ws.Range("Q3").Value = 606.3267665309306
$ws.Range("R3").Value = 5456.940898778375
$ws.Range("S3").Value = 0.8617236586166084
$ws.Range("T3").Value = 0.8708473271753914

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.447252
$ws.Range("H4").Value = 61.341756
$ws.Range("I4").Value = 0.8699145605694745
$ws.Range("J4").Value = 0.8770588936480435
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04031433333333333
$ws.Range("N4").Value = 0.120943
$ws.Range("O4").Value = 0.001346725574891351
$ws.Range("P4").Value = 0.001349898011227873
$ws.Range("Q4").Value = 0.8243173328786666
$ws.Range("R4").Value = 7.418855995907999
$ws.Range("S4").Value = 0.001171536186689282
$ws.Range("T4").Value = 0.001183940056265212

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.447252
$ws.Range("H5").Value = 61.341756
$ws.Range("I5").Value = 0.8699145605694745
$ws.Range("J5").Value = 0.8770588936480435
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.211054
$ws.Range("N5").Value = 0.422108
$ws.Range("O5").Value = 0.007050391163187267
$ws.Range("P5").Value = 0.00471133302236074
$ws.Range("Q5").Value = 4.315474323608
$ws.Range("R5").Value = 25.892845941648
$ws.Range("S5").Value = 0.006133237930566957
$ws.Range("T5").Value = 0.004132116528199203

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.483247333333333
$ws.Range("H6").Value = 7.449742
$ws.Range("I6").Value = 0.105648084777455
$ws.Range("J6").Value = 0.1065157390747562
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.030493
$ws.Range("N6").Value = 0.091479
$ws.Range("O6").Value = 0.001018637778668347
$ws.Range("P6").Value = 0.001021037349570579
$ws.Range("Q6").Value = 0.07572166093533333
$ws.Range("R6").Value = 0.681494948418
$ws.Range("S6").Value = 0.000107617130398272
$ws.Range("T6").Value = 0.0001087565479124405

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.483247333333333
$ws.Range("H7").Value = 7.449742
$ws.Range("I7").Value = 0.105648084777455
$ws.Range("J7").Value = 0.1065157390747562
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 29.65321533333333
$ws.Range("N7").Value = 88.95964599999999
$ws.Range("O7").Value = 0.990584245483253
$ws.Range("P7").Value = 0.9929177316168408
$ws.Range("Q7").Value = 73.6362679012591
$ws.Range("R7").Value = 662.7264111113319
$ws.Range("S7").Value = 0.104653328346026
$ws.Range("T7").Value = 0.1057613660235983

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.483247333333333
$ws.Range("H8").Value = 7.449742
$ws.Range("I8").Value = 0.105648084777455
$ws.Range("J8").Value = 0.1065157390747562
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04031433333333333
$ws.Range("N8").Value = 0.120943
$ws.Range("O8").Value = 0.001346725574891351
$ws.Range("P8").Value = 0.001349898011227873
$ws.Range("Q8").Value = 0.1001104607451111
$ws.Range("R8").Value = 0.9009941467059999
$ws.Range("S8").Value = 0.0001422789777080883
$ws.Range("T8").Value = 0.0001437853843414805

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.483247333333333
$ws.Range("H9").Value = 7.449742
$ws.Range("I9").Value = 0.105648084777455
$ws.Range("J9").Value = 0.1065157390747562
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.211054
$ws.Range("N9").Value = 0.422108
$ws.Range("O9").Value = 0.007050391163187267
$ws.Range("P9").Value = 0.00471133302236074
$ws.Range("Q9").Value = 0.5240992826893334
$ws.Range("R9").Value = 3.144595696136
$ws.Range("S9").Value = 0.0007448603233226279
$ws.Range("T9").Value = 0.0005018311189040592

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5743975
$ws.Range("H10").Value = 1.148795
$ws.Range("I10").Value = 0.02443735465307048
$ws.Range("J10").Value = 0.01642536727720028
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.030493
$ws.Range("N10").Value = 0.091479
$ws.Range("O10").Value = 0.001018637778668347
$ws.Range("P10").Value = 0.001021037349570579
$ws.Range("Q10").Value = 0.0175151029675
$ws.Range("R10").Value = 0.105090617805
$ws.Range("S10").Value = 0.00002489281266033432
$ws.Range("T10").Value = 0.0000167709134704359

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Efna1"
$ws.Range("C11").Value = "Epha3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5743975
$ws.Range("H11").Value = 1.148795
$ws.Range("I11").Value = 0.02443735465307048
$ws.Range("J11").Value = 0.01642536727720028
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 29.65321533333333
$ws.Range("N11").Value = 88.95964599999999
$ws.Range("O11").Value = 0.990584245483253
$ws.Range("P11").Value = 0.9929177316168408
$ws.Range("Q11").Value = 17.03273275442833
$ws.Range("R11").Value = 102.19639652657
$ws.Range("S11").Value = 0.02420725852061849
$ws.Range("T11").Value = 0.01630903841785119

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Efna1"
$ws.Range("C12").Value = "Epha3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5743975
$ws.Range("H12").Value = 1.148795
$ws.Range("I12").Value = 0.02443735465307048
$ws.Range("J12").Value = 0.01642536727720028
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.04031433333333333
$ws.Range("N12").Value = 0.120943
$ws.Range("O12").Value = 0.001346725574891351
$ws.Range("P12").Value = 0.001349898011227873
$ws.Range("Q12").Value = 0.02315645228083333
$ws.Range("R12").Value = 0.138938713685
$ws.Range("S12").Value = 0.00003291041049398018
$ws.Range("T12").Value = 0.00002217257062118004

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Efna1"
$ws.Range("C13").Value = "Epha3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5743975
$ws.Range("H13").Value = 1.148795
$ws.Range("I13").Value = 0.02443735465307048
$ws.Range("J13").Value = 0.01642536727720028
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.211054
$ws.Range("N13").Value = 0.422108
$ws.Range("O13").Value = 0.007050391163187267
$ws.Range("P13").Value = 0.00471133302236074
$ws.Range("Q13").Value = 0.121228889965
$ws.Range("R13").Value = 0.4849155598600001
$ws.Range("S13").Value = 0.0001722929092976814
$ws.Range("T13").Value = 0.00007738537525747721
